$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - column F updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 0
$ws1.Range("F3").Value = 635
$ws1.Range("F4").Value = 0
$ws1.Range("F5").Value = 4946
$ws1.Range("F6").Value = 516
$ws1.Range("F7").Value = 9118
$ws1.Range("F9").Value = 512
$ws1.Range("F12").Value = 68

# Sheet "演出" (shows) - column F updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 0
$ws2.Range("F3").Value = 0
$ws2.Range("F4").Value = 8
$ws2.Range("F5").Value = 6
$ws2.Range("F6").Value = 3

# Sheet "全部类型" (all types) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1151
$ws4.Range("F3").Value = 635
$ws4.Range("F4").Value = 344
$ws4.Range("F5").Value = 14
$ws4.Range("F6").Value = 20
$ws4.Range("F7").Value = 4946
$ws4.Range("F9").Value = 8
$ws4.Range("F10").Value = 0
$ws4.Range("F12").Value = 512
$ws4.Range("F14").Value = 6
$ws4.Range("F16").Value = 0
